$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 70: nbr2
$ws.Range("A70").Value = "satellite"
$ws.Range("B70").Value = "landsat"
$ws.Range("C70").Value = "AusCover"
$ws.Range("D70").Value = "surface reflectance"
$ws.Range("E70").Value = "nbr2"
$ws.Range("F70").Value = 1
$ws.Range("G70").Value = "mean"
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = "each"
$ws.Range("J70").Value = 1989
$ws.Range("J70").NumberFormat = "0"
$ws.Range("K70").Value = 10000
$ws.Range("L70").Value = "units"
$ws.Range("N70").Value = "(nir - swir1)/(nir + swir1)"
$ws.Range("O70").Value = "modifies the Normalized Burn Ratio (NBR) to highlight water sensitivity in vegetation and may be useful in post-fire recovery studies"
$ws.Range("M70").Value = "normalised burn ratio 2"

# Row 71: savi
$ws.Range("A71").Value = "satellite"
$ws.Range("B71").Value = "landsat"
$ws.Range("C71").Value = "AusCover"
$ws.Range("D71").Value = "surface reflectance"
$ws.Range("E71").Value = "savi"
$ws.Range("F71").Value = 1
$ws.Range("G71").Value = "mean"
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = "each"
$ws.Range("J71").Value = 1989
$ws.Range("J71").NumberFormat = "0"
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = "units"
$ws.Range("M71").Value = "soil adjusted vegetation index"
$ws.Range("N71").Value = "((nir - red)/(nir + red + L)) * (1 + L)"
$ws.Range("O71").Value = "used to correct Normalized Difference Vegetation Index (NDVI) for the influence of soil brightness in areas where vegetative cover is low"

# Update the view's active selection to match the authored state
$ws.Range("G67").Select()
